# Dev IV Solo Project Rubric - grading update
# "CommandList unfinished, Drawinstance is kind of buggy"
#
# This updates the milestone achievement grid (columns E = milestone level
# achieved, F = milestone completed flag) for several rubric rows, adds a
# missing source citation, and moves the active cell selection. All of the
# summary/point totals (G, H, I, J, K, L columns) are formula driven and will
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22: "Apply multitexturing or advanced texturing techniques..." now
# marked complete at Milestone II.
$ws.Range("E22").Value2 = "II"
$ws.Range("F22").Value2 = "X"

# Row 38: "Drawing procedurally created 3D line mesh..." - unmark Milestone
# III assignment (left blank again).
$ws.Range("E38").Value2 = ""

# Row 40: "Use DrawInstanced() or DrawIndexedInstanced()..." now targeted at
# Milestone III (per commit message, DrawInstanced support is being worked on).
$ws.Range("E40").Value2 = "III"

# Row 65: "Camera tracks moving object in scene using Look-At/Turn-To..."
# moved from Milestone II to Milestone III.
$ws.Range("E65").Value2 = "III"

# Row 67: "Multiple viewport scenes such as minimaps or splitscreen..."
# moved from Milestone II to Milestone III.
$ws.Range("E67").Value2 = "III"

# Row 85: "Island (Tropical) Theme..." now marked as completed.
$ws.Range("F85").Value2 = "X"

# Row 96: add the missing skybox source citation.
$ws.Range("A96").Value2 = "space Skybox -- Created by 'amethyst7' aka Chris Matz. WEB site: http://amethyst7.gotdoofed.com"

# Recalculate all formulas (totals, grades, carry-over, etc.)
$excel.CalculateFullRebuild()

# Move the selected/active cell as it was left in the saved workbook.
$ws.Range("E60").Select()
